$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "ID" column
$ws.Range("C1").Value = "ID"

# Fill in student IDs for rows 2-16 (100001 .. 100015)
$id = 100001
for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = $id
    $id = $id + 1
}

# Update the selected cell to E7
$null = $ws.Range("E7").Select()
